$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D:E (price / volume) to Text format before writing values so
# that numeric-looking strings (e.g. "42.870.58", "94.16", "  -1.63%  ")
# are preserved as text rather than being auto-converted to numbers, matching
# the original inlineStr cell type used throughout the sheet.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '42.870.58'
$ws.Range("E2").Value = '  -1.63%  '
$ws.Range("D3").Value = '2.556.73'
$ws.Range("E3").Value = '  -2.69%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = '301.22'
$ws.Range("E5").Value = '  -0.04%  '
$ws.Range("D6").Value = '94.16'
$ws.Range("E6").Value = '  -2.19%  '
$ws.Range("E7").Value = '  -1.54%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("D9").Value = '0.544'
$ws.Range("E9").Value = '  -2.19%  '
$ws.Range("D10").Value = '36.31'
$ws.Range("E10").Value = '  -2.14%  '
$ws.Range("E11").Value = '  -0.61%  '
$ws.Range("D12").Value = '7.73'
$ws.Range("E12").Value = '  -1.08%  '
$ws.Range("E13").Value = '  +7.07%  '
$ws.Range("D14").Value = '2.952.08'
$ws.Range("E14").Value = '  -2.63%  '
$ws.Range("D15").Value = '2.577.57'
$ws.Range("E15").Value = '  -1.65%  '
$ws.Range("D16").Value = '0.876'
$ws.Range("E16").Value = '  -1.58%  '
$ws.Range("D17").Value = '14.16'
$ws.Range("E17").Value = '  -1.84%  '
$ws.Range("D18").Value = '42.913.12'
$ws.Range("E18").Value = '  -1.50%  '
$ws.Range("D19").Value = '0.0₃0985'
$ws.Range("E19").Value = '  +0.92%  '
$ws.Range("D20").Value = '12.66'
$ws.Range("E20").Value = '  +1.53%  '
$ws.Range("D21").Value = '6.53'
$ws.Range("E21").Value = '  -1.98%  '
$ws.Range("D22").Value = '71.43'
$ws.Range("E22").Value = '  -3.14%  '
$ws.Range("D23").Value = '252.18'
$ws.Range("E23").Value = '  -5.94%  '
$ws.Range("D24").Value = '2.94'
$ws.Range("E24").Value = '  +0.13%  '
$ws.Range("D25").Value = '2.11'
$ws.Range("E25").Value = '  -5.15%  '
$ws.Range("D26").Value = '28.67'
$ws.Range("E26").Value = '  -3.45%  '
$ws.Range("E27").Value = '  -0.18%  '
$ws.Range("D28").Value = '10.20'
$ws.Range("E28").Value = '  -0.44%  '
$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").Value = '2.14'
$ws.Range("E29").Value = '  -4.34%  '
$ws.Range("B30").Value = 'InjectiveProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D30").Value = '36.76'
$ws.Range("E30").Value = '  -3.09%  '
$ws.Range("D31").Value = '6.05'
$ws.Range("E31").Value = '  +0.04%  '
$ws.Range("D32").Value = '153.81'
$ws.Range("E32").Value = '  +0.87%  '
$ws.Range("D33").Value = '2.74'
$ws.Range("E33").Value = '  -1.75%  '
$ws.Range("D34").Value = '3.36'
$ws.Range("E34").Value = '  -7.48%  '
$ws.Range("D35").Value = '2.13'
$ws.Range("E35").Value = '  -5.31%  '
$ws.Range("D36").Value = '0.0796'
$ws.Range("E36").Value = '  -2.00%  '
$ws.Range("D37").Value = '0.113'
$ws.Range("E37").Value = '  -3.76%  '
$ws.Range("D38").Value = '17.81'
$ws.Range("E38").Value = '  +10.07%  '
$ws.Range("E39").Value = '  -0.63%  '
$ws.Range("D40").Value = '23.33'
$ws.Range("E40").Value = '  -6.87%  '
$ws.Range("D41").Value = '2.11'
$ws.Range("E41").Value = '  +32.14%  '
$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D42").Value = '0.0309'
$ws.Range("E42").Value = '  -1.71%  '
$ws.Range("B43").Value = 'NEARProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D43").Value = '3.39'
$ws.Range("E43").Value = '  -3.12%  '
$ws.Range("D44").Value = '3.85'
$ws.Range("E44").Value = '  -0.21%  '
$ws.Range("D45").Value = '2.098.82'
$ws.Range("E45").Value = '  +0.28%  '
$ws.Range("E46").Value = '  +0.18%  '
$ws.Range("D47").Value = '9.24'
$ws.Range("E47").Value = '  +1.64%  '
$ws.Range("D48").Value = '84.68'
$ws.Range("E48").Value = '  -5.11%  '
$ws.Range("B49").Value = 'RocketPoolETH'
$ws.Range("C49").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D49").Value = '2.808.49'
$ws.Range("E49").Value = '  -2.94%  '
$ws.Range("B50").Value = 'Aave'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D50").Value = '105.56'
$ws.Range("E50").Value = '  -1.17%  '
$ws.Range("D51").Value = '73.99'
$ws.Range("E51").Value = '  +6.96%  '

# Remove the temporary Text number-format we applied above so the cells end
# up with no explicit style (matching the original workbook, where these
# data cells carry no s= attribute).
$ws.Range("D2:E51").ClearFormats()
